$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cumulative AgTests (F) / AgPosit (G) figures, and fill in the
# previously-missing F/G values for row 727, per the latest data refresh.
$ws.Cells.Item(484, 6).Value = 8452
$ws.Cells.Item(546, 6).Value = 4055
$ws.Cells.Item(576, 6).Value = 29570
$ws.Cells.Item(637, 6).Value = 43804
$ws.Cells.Item(671, 6).Value = 32667
$ws.Cells.Item(678, 6).Value = 33874
$ws.Cells.Item(679, 6).Value = 29485
$ws.Cells.Item(680, 6).Value = 28476
$ws.Cells.Item(681, 6).Value = 26441
$ws.Cells.Item(683, 6).Value = 24290
$ws.Cells.Item(684, 6).Value = 57276
$ws.Cells.Item(684, 7).Value = 1214
$ws.Cells.Item(685, 6).Value = 34492
$ws.Cells.Item(686, 6).Value = 34457
$ws.Cells.Item(687, 6).Value = 31504
$ws.Cells.Item(688, 6).Value = 32191
$ws.Cells.Item(690, 6).Value = 27889
$ws.Cells.Item(690, 7).Value = 1550
$ws.Cells.Item(691, 6).Value = 62643
$ws.Cells.Item(692, 6).Value = 41690
$ws.Cells.Item(693, 6).Value = 39621
$ws.Cells.Item(693, 7).Value = 2734
$ws.Cells.Item(694, 6).Value = 37639
$ws.Cells.Item(694, 7).Value = 2778
$ws.Cells.Item(695, 6).Value = 37308
$ws.Cells.Item(695, 7).Value = 3138
$ws.Cells.Item(697, 6).Value = 28959
$ws.Cells.Item(697, 7).Value = 3042
$ws.Cells.Item(698, 6).Value = 70826
$ws.Cells.Item(698, 7).Value = 5825
$ws.Cells.Item(699, 6).Value = 43562
$ws.Cells.Item(699, 7).Value = 4306
$ws.Cells.Item(700, 6).Value = 43749
$ws.Cells.Item(700, 7).Value = 4319
$ws.Cells.Item(701, 6).Value = 41854
$ws.Cells.Item(701, 7).Value = 3853
$ws.Cells.Item(708, 6).Value = 35514
$ws.Cells.Item(709, 6).Value = 32330
$ws.Cells.Item(712, 6).Value = 51253
$ws.Cells.Item(713, 6).Value = 37085
$ws.Cells.Item(715, 6).Value = 31640
$ws.Cells.Item(716, 6).Value = 29552
$ws.Cells.Item(718, 6).Value = 16817
$ws.Cells.Item(718, 7).Value = 2793
$ws.Cells.Item(719, 6).Value = 43313
$ws.Cells.Item(719, 7).Value = 5144
$ws.Cells.Item(720, 6).Value = 30911
$ws.Cells.Item(720, 7).Value = 3479
$ws.Cells.Item(721, 6).Value = 27706
$ws.Cells.Item(721, 7).Value = 3100
$ws.Cells.Item(722, 6).Value = 27685
$ws.Cells.Item(722, 7).Value = 2835
$ws.Cells.Item(723, 6).Value = 22039
$ws.Cells.Item(723, 7).Value = 2709
$ws.Cells.Item(724, 6).Value = 9257
$ws.Cells.Item(724, 7).Value = 1476
$ws.Cells.Item(725, 6).Value = 12416
$ws.Cells.Item(725, 7).Value = 2020
$ws.Cells.Item(726, 6).Value = 34506
$ws.Cells.Item(726, 7).Value = 3995

# Row 727 previously had no AgTests/AgPosit values; fill them in
$ws.Cells.Item(727, 6).Value = 23760
$ws.Cells.Item(727, 7).Value = 2672

# Append new row 728 with data for 2022-03-02 (serial date 44622)
$ws.Cells.Item(728, 1).Value = 44622
$ws.Cells.Item(728, 2).Value = 1482354
$ws.Cells.Item(728, 3).Value = 20128
$ws.Cells.Item(728, 4).Value = 11462
$ws.Cells.Item(728, 5).Value = 18611
$ws.Cells.Item(728, 6).Value = 17297
$ws.Cells.Item(728, 7).Value = 1784

